$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row values (B1:E1)
$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

# Update row 2 values (B2:E2)
$ws.Range("B2").Value = 4.7037690973825077
$ws.Range("C2").Value = 1.2637975731455251
$ws.Range("D2").Value = 0.87508580423249305
$ws.Range("E2").Value = 0.30266655436741863

# Update row 3 values (B3:E3)
$ws.Range("B3").Value = 4.2882762550519846
$ws.Range("C3").Value = 5.8867350382756936
$ws.Range("D3").Value = 5.8290072667893034
$ws.Range("E3").Value = -1.725174307024737

# Update the selection to match the new range B1:E3
$ws.Range("B1:E3").Select()
